# Update Name of Algo
# Refresh the 61 KNN-imputed B/D/E values in Sheet1 (rows 3-102) that
# changed when the algorithm result was regenerated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = 16.349
$ws.Range("E4").Value = 16.554
$ws.Range("D7").Value = -7.846000000000001
$ws.Range("B8").Value = 6.223000000000001
$ws.Range("B10").Value = 5.816
$ws.Range("E11").Value = 17.024
$ws.Range("B12").Value = 5.724
$ws.Range("D14").Value = -7.933
$ws.Range("E14").Value = 16.816
$ws.Range("D15").Value = -8.196000000000002
$ws.Range("B18").Value = 4.972
$ws.Range("D18").Value = -8.32
$ws.Range("E18").Value = 16.583
$ws.Range("E19").Value = 16.452
$ws.Range("D20").Value = -7.579000000000001
$ws.Range("E21").Value = 16.421
$ws.Range("B25").Value = 6.074
$ws.Range("E27").Value = 16.473
$ws.Range("D29").Value = -7.394
$ws.Range("D30").Value = -7.048999999999999
$ws.Range("D31").Value = -7.613000000000001
$ws.Range("E31").Value = 16.822
$ws.Range("D35").Value = -7.737
$ws.Range("B37").Value = 8.56
$ws.Range("E38").Value = 16.731
$ws.Range("D40").Value = -7.587999999999999
$ws.Range("E42").Value = 16.425
$ws.Range("D44").Value = -7.447
$ws.Range("E44").Value = 16.992
$ws.Range("E47").Value = 16.442
$ws.Range("D50").Value = -8.125999999999999
$ws.Range("D54").Value = -8.010999999999999
$ws.Range("B55").Value = 4.852
$ws.Range("E56").Value = 16.313
$ws.Range("E58").Value = 16.387
$ws.Range("E65").Value = 17.023
$ws.Range("B68").Value = 5.788
$ws.Range("D68").Value = -6.738
$ws.Range("E73").Value = 16.452
$ws.Range("D76").Value = -7.444
$ws.Range("B77").Value = 5.474
$ws.Range("B78").Value = 7.538000000000001
$ws.Range("B79").Value = 6.375
$ws.Range("B80").Value = 8.183
$ws.Range("B81").Value = 5.879
$ws.Range("B82").Value = 5.659000000000001
$ws.Range("B84").Value = 5.921
$ws.Range("D87").Value = -8.297000000000001
$ws.Range("D88").Value = -8.056000000000001
$ws.Range("E90").Value = 16.49
$ws.Range("D92").Value = -7.017
$ws.Range("E92").Value = 17.237
$ws.Range("E94").Value = 17.699
$ws.Range("E95").Value = 17.136
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.244000000000002
$ws.Range("B101").Value = 6.813999999999998
$ws.Range("D101").Value = -7.756
$ws.Range("E101").Value = 16.383
$ws.Range("B102").Value = 8.106
$ws.Range("D102").Value = -7.929
